# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# For each row, Price (D) and/or Volume(1h) (E) text is replaced with the new scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "22.373.55"
$ws.Range("E2").Value2 = "  -4.55%  "
$ws.Range("D3").Value2 = "1.564.48"
$ws.Range("E3").Value2 = "  -5.00%  "
$ws.Range("D4").Value2 = "`'1.001"
$ws.Range("E4").Value2 = "  +0.21%  "
$ws.Range("E5").Value2 = "  +0.26%  "
$ws.Range("D6").Value2 = "`'289.76"
$ws.Range("E6").Value2 = "  -3.28%  "
$ws.Range("D7").Value2 = "`'0.3706"
$ws.Range("E7").Value2 = "  -2.42%  "
$ws.Range("D8").Value2 = "`'49.24"
$ws.Range("E8").Value2 = "  -1.93%  "
$ws.Range("D9").Value2 = "`'0.3385"
$ws.Range("E9").Value2 = "  -3.21%  "
$ws.Range("D10").Value2 = "`'1.165"
$ws.Range("E10").Value2 = "  -4.36%  "
$ws.Range("D11").Value2 = "`'0.07645"
$ws.Range("E11").Value2 = "  -5.26%  "
$ws.Range("E12").Value2 = "  +0.35%  "
$ws.Range("D13").Value2 = "`'21.42"
$ws.Range("E13").Value2 = "  -2.84%  "
$ws.Range("D14").Value2 = "`'6.039"
$ws.Range("E14").Value2 = "  -4.49%  "
$ws.Range("D15").Value2 = "`'6.920"
$ws.Range("E15").Value2 = "  -4.90%  "
$ws.Range("D16").Value2 = "1.570.99"
$ws.Range("E16").Value2 = "  -3.67%  "
$ws.Range("D17").Value2 = "`'0.00001127"
$ws.Range("E17").Value2 = "  -7.18%  "
$ws.Range("D18").Value2 = "`'90.16"
$ws.Range("E18").Value2 = "  -4.94%  "
$ws.Range("D19").Value2 = "`'0.06728"
$ws.Range("E19").Value2 = "  -3.42%  "
$ws.Range("E20").Value2 = "  +0.32%  "
$ws.Range("D21").Value2 = "`'6.239"
$ws.Range("E21").Value2 = "  -5.94%  "
$ws.Range("D22").Value2 = "`'16.53"
$ws.Range("E22").Value2 = "  -4.93%  "
$ws.Range("D23").Value2 = "`'0.5280"
$ws.Range("E23").Value2 = "  -7.94%  "
$ws.Range("D24").Value2 = "`'12.02"
$ws.Range("E24").Value2 = "  -3.11%  "
$ws.Range("D25").Value2 = "22.359.76"
$ws.Range("E25").Value2 = "  -4.67%  "
$ws.Range("D26").Value2 = "`'2.369"
$ws.Range("E26").Value2 = "  -2.70%  "
$ws.Range("D27").Value2 = "`'2.804"
$ws.Range("E27").Value2 = "  -5.50%  "
$ws.Range("D28").Value2 = "`'20.15"
$ws.Range("E28").Value2 = "  -4.21%  "
$ws.Range("D29").Value2 = "`'145.24"
$ws.Range("E29").Value2 = "  -3.26%  "
$ws.Range("D30").Value2 = "`'4.984"
$ws.Range("E30").Value2 = "  -3.81%  "
$ws.Range("D31").Value2 = "`'125.58"
$ws.Range("E31").Value2 = "  -4.75%  "
$ws.Range("D32").Value2 = "1.724.06"
$ws.Range("E32").Value2 = "  -5.13%  "
$ws.Range("D33").Value2 = "`'6.220"
$ws.Range("E33").Value2 = "  -9.24%  "
$ws.Range("D34").Value2 = "`'2.007"
$ws.Range("E34").Value2 = "  -5.71%  "
$ws.Range("D35").Value2 = "`'1.001"
$ws.Range("E35").Value2 = "  +1.11%  "
$ws.Range("D36").Value2 = "`'10.10"
$ws.Range("E36").Value2 = "  -10.33%  "
$ws.Range("D37").Value2 = "`'0.08439"
$ws.Range("E37").Value2 = "  -3.97%  "
$ws.Range("D38").Value2 = "`'0.02533"
$ws.Range("E38").Value2 = "  -5.70%  "
$ws.Range("D39").Value2 = "`'0.2323"
$ws.Range("E39").Value2 = "  -4.16%  "
$ws.Range("D40").Value2 = "`'5.525"
$ws.Range("E40").Value2 = "  -6.21%  "
$ws.Range("D41").Value2 = "`'0.06418"
$ws.Range("E41").Value2 = "  -6.11%  "
$ws.Range("D42").Value2 = "`'1.295"
$ws.Range("E42").Value2 = "  +0.61%  "
$ws.Range("D43").Value2 = "`'11.71"
$ws.Range("E43").Value2 = "  -8.39%  "
$ws.Range("D44").Value2 = "`'0.6334"
$ws.Range("E44").Value2 = "  -7.31%  "
$ws.Range("D45").Value2 = "`'14.25"
$ws.Range("E45").Value2 = "  -7.75%  "
$ws.Range("E46").Value2 = "  +0.32%  "
$ws.Range("D47").Value2 = "`'0.5979"
$ws.Range("E47").Value2 = "  -5.90%  "
$ws.Range("E48").Value2 = "  -4.03%  "
$ws.Range("D49").Value2 = "`'2.097"
$ws.Range("E49").Value2 = "  -6.44%  "
$ws.Range("D50").Value2 = "`'1.264"
$ws.Range("E50").Value2 = "  +2.96%  "
$ws.Range("D51").Value2 = "`'124.33"
$ws.Range("E51").Value2 = "  -2.27%  "
